$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Variables Listing File name from v4 to v5
$ws.Range("B13").Value = "trend_report_variables_v5.xlsx"

# Update the Variables List Indices lower right cell reference
$ws.Range("D15").Value = "E238"

# Update the active selection to match the saved view state
$ws.Range("D16").Select()
